$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1872146118721461
$ws.Range("C2").Value = 0.5958904109589042
$ws.Range("J2").Value = 0.01598173515981735
$ws.Range("P2").Value = 0.1324200913242009
$ws.Range("S2").Value = 0.0684931506849315
$ws.Range("B3").Value = 0.003731343283582089
$ws.Range("C3").Value = 0.01492537313432836
$ws.Range("J3").Value = 0.02611940298507463
$ws.Range("P3").Value = 0.7649253731343284
$ws.Range("S3").Value = 0.1902985074626866
$ws.Range("J4").Value = 0.1020408163265306
$ws.Range("P4").Value = 0.673469387755102
$ws.Range("S4").Value = 0.2244897959183673
$ws.Range("B6").Value = 0.0641025641025641
$ws.Range("D6").Value = 0.02136752136752137
$ws.Range("E6").Value = 0.004273504273504274
$ws.Range("F6").Value = 0.05982905982905983
$ws.Range("J6").Value = 0.2521367521367521
$ws.Range("O6").Value = 0.01282051282051282
$ws.Range("Q6").Value = 0.1495726495726496
$ws.Range("R6").Value = 0.08974358974358974
$ws.Range("B7").Value = 0.1490909090909091
$ws.Range("D7").Value = 0.02909090909090909
$ws.Range("F7").Value = 0.08
$ws.Range("J7").Value = 0.1054545454545455
$ws.Range("O7").Value = 0.01090909090909091
$ws.Range("Q7").Value = 0.1527272727272727
$ws.Range("R7").Value = 0.06545454545454546
$ws.Range("S7").Value = 0.4072727272727273
$ws.Range("B8").Value = 0.1266094420600858
$ws.Range("D8").Value = 0.01716738197424893
$ws.Range("F8").Value = 0.06008583690987124
$ws.Range("J8").Value = 0.1266094420600858
$ws.Range("O8").Value = 0.01931330472103004
$ws.Range("Q8").Value = 0.1437768240343348
$ws.Range("R8").Value = 0.09012875536480687
$ws.Range("S8").Value = 0.4163090128755365
$ws.Range("B9").Value = 0.1171171171171171
$ws.Range("D9").Value = 0.02702702702702703
$ws.Range("F9").Value = 0.07207207207207207
$ws.Range("J9").Value = 0.1801801801801802
$ws.Range("O9").Value = 0.05405405405405406
$ws.Range("Q9").Value = 0.1351351351351351
$ws.Range("R9").Value = 0.07207207207207207
$ws.Range("S9").Value = 0.3423423423423423
$ws.Range("B10").Value = 0.1482461945731304
$ws.Range("D10").Value = 0.01786896095301125
$ws.Range("F10").Value = 0.06022501654533421
$ws.Range("J10").Value = 0.1317008603573792
$ws.Range("O10").Value = 0.01389808074123097
$ws.Range("Q10").Value = 0.1833223031105228
$ws.Range("R10").Value = 0.08802117802779616
$ws.Range("S10").Value = 0.356717405691595
$ws.Range("G11").Value = 0.1363636363636364
$ws.Range("J11").Value = 0.08636363636363636
$ws.Range("K11").Value = 0.1977272727272727
$ws.Range("L11").Value = 0.5613636363636364
$ws.Range("S11").Value = 0.01818181818181818
$ws.Range("G12").Value = 0.7261904761904762
$ws.Range("J12").Value = 0.2063492063492063
$ws.Range("K12").Value = 0.0119047619047619
$ws.Range("L12").Value = 0.0119047619047619
$ws.Range("S12").Value = 0.04365079365079365
$ws.Range("G13").Value = 0.6507936507936508
$ws.Range("J13").Value = 0.2698412698412698
$ws.Range("S13").Value = 0.07936507936507936
$ws.Range("F15").Value = 0.004950495049504951
$ws.Range("H15").Value = 0.1435643564356436
$ws.Range("I15").Value = 0.04455445544554455
$ws.Range("J15").Value = 0.3663366336633663
$ws.Range("K15").Value = 0.03465346534653466
$ws.Range("M15").Value = 0.01485148514851485
$ws.Range("N15").Value = 0.004950495049504951
$ws.Range("O15").Value = 0.04950495049504951
$ws.Range("S15").Value = 0.3366336633663367
$ws.Range("F16").Value = 0.02068965517241379
$ws.Range("H16").Value = 0.2068965517241379
$ws.Range("I16").Value = 0.05517241379310345
$ws.Range("J16").Value = 0.4137931034482759
$ws.Range("K16").Value = 0.1310344827586207
$ws.Range("M16").Value = 0.01379310344827586
$ws.Range("O16").Value = 0.03448275862068965
$ws.Range("S16").Value = 0.1241379310344828
$ws.Range("F17").Value = 0.01392111368909513
$ws.Range("H17").Value = 0.1647331786542924
$ws.Range("I17").Value = 0.04408352668213457
$ws.Range("J17").Value = 0.4129930394431555
$ws.Range("K17").Value = 0.160092807424594
$ws.Range("M17").Value = 0.02088167053364269
$ws.Range("O17").Value = 0.06032482598607888
$ws.Range("S17").Value = 0.122969837587007
$ws.Range("F18").Value = 0.01809954751131222
$ws.Range("H18").Value = 0.1493212669683258
$ws.Range("I18").Value = 0.05882352941176471
$ws.Range("J18").Value = 0.4208144796380091
$ws.Range("K18").Value = 0.1493212669683258
$ws.Range("M18").Value = 0.03167420814479638
$ws.Range("O18").Value = 0.07239819004524888
$ws.Range("S18").Value = 0.09954751131221719
$ws.Range("F19").Value = 0.015625
$ws.Range("H19").Value = 0.1960227272727273
$ws.Range("I19").Value = 0.03835227272727273
$ws.Range("J19").Value = 0.4055397727272727
$ws.Range("K19").Value = 0.1420454545454546
$ws.Range("M19").Value = 0.02840909090909091
$ws.Range("N19").Value = 0.001420454545454545
$ws.Range("O19").Value = 0.05397727272727273
$ws.Range("S19").Value = 0.1186079545454545
